$d = $word.ActiveDocument

function MergeRuns($searchText) {
    # Find the (unique) occurrence of $searchText and replace it with the
    # exact same text. Because the text spans two (or more) adjacent runs
    # that share identical run formatting, Word coalesces them into a
    # single run as a side effect of the replace.
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $searchText, 2)
    if (-not $found) {
        Write-Host "WARNING: not found:" $searchText
    }
    return $found
}

# 1. " problem 3-122 (see Appendix). To determine the " + "maximum force, "
MergeRuns(" problem 3-122 (see Appendix). To determine the maximum force, ")

# 2. ", that can safely be applied to the eye" + "bolt, apply a nominal force, "
MergeRuns(", that can safely be applied to the eyebolt, apply a nominal force, ")

# 3. ", that can safely be applied to the eyebolt" + " is"
MergeRuns(", that can safely be applied to the eyebolt is")

# 5. "and the equation given above, calculate " + "the maximum force, "
MergeRuns("and the equation given above, calculate the maximum force, ")

# 6. "O" + "n the Solution Summary Form" + " (attached)...the " + "von Mises stress, "
MergeRuns("On the Solution Summary Form (attached), record the global element size used in your final mesh, the element type (tetrahedral or hexahedral), interpolation type (linear or quadratic), the von Mises stress, ")

# 7. "=1 N" + " and " + "the maximum force, "
MergeRuns("=1 N and the maximum force, ")

# 8. ", that can safely be applied to the eyebolt" + "." -- this phrase occurs
# twice in the document; only the *second* occurrence (the one following
# "=1 N and the maximum force, ...") is merged by the target edit, so we
# must locate it positionally rather than simply searching document-wide.
$probe = "that can safely be applied to the eyebolt."
$r1 = $d.Content
$r1.Find.Execute($probe, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2 = $d.Range($r1.End, $d.Content.End)
$r2.Find.Execute($probe, $true, $false, $false, $false, $false, $true, 1, $false, $probe, 2) | Out-Null

# 9. "at the inner radius of section A-A" + " for an applied force, " (inside table)
MergeRuns("at the inner radius of section A-A for an applied force, ")

# 10. "=1 N" + " " (inside table) -- careful: the surrounding "Fnom" and
# "(MPa)" runs use different formatting (italic / different eastAsia font),
# so we must only touch the "=1 N " span itself, not the whole phrase.
$probe10 = "Fnom=1 N (MPa)"
$rng10 = $d.Content
$rng10.Find.Execute($probe10, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sub10 = $d.Range($rng10.Start + 4, $rng10.Start + 9)
$sub10txt = $sub10.Text
$sub10.Find.Execute($sub10txt, $true, $false, $false, $false, $false, $true, 1, $false, $sub10txt, 2) | Out-Null

# 4. Remove the old (stray) _GoBack bookmark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 11. Split " under Computer Assignment #5" into " under Computer Assignment #"
# + "4", then re-insert the _GoBack bookmark right after the new "4" run.
$rngCA = $d.Content
$rngCA.Find.Execute("under Computer Assignment #5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $rngCA.End

$digit = $d.Range($endPos - 1, $endPos)
$digit.Text = "4"

# Force a genuine run split for the freshly typed "4" (toggling Bold on/off
# leaves formatting unchanged but breaks it out of the surrounding run).
$fourRange = $d.Range($endPos - 1, $endPos)
$fourRange.Bold = 1
$fourRange.Bold = 0

# Re-insert the _GoBack bookmark immediately after the "4".
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
